$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-07-31 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-01 Friday", 2)

$d.Content.Find.Execute("36÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "25÷9=", 2)
$d.Content.Find.Execute("74÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "39÷9=", 2)
$d.Content.Find.Execute("19÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "40÷9=", 2)
$d.Content.Find.Execute("58÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "24÷9=", 2)
$d.Content.Find.Execute("91÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "34÷5=", 2)

$d.Content.Find.Execute("11÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "51÷8=", 2)
$d.Content.Find.Execute("29÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "72÷6=", 2)
$d.Content.Find.Execute("80÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "66÷2=", 2)
$d.Content.Find.Execute("16÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "98÷9=", 2)
$d.Content.Find.Execute("11÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "95÷9=", 2)

$d.Content.Find.Execute("86÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "45÷7=", 2)
$d.Content.Find.Execute("49÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "84÷2=", 2)
$d.Content.Find.Execute("93÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "91÷3=", 2)
$d.Content.Find.Execute("43÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "85÷5=", 2)
$d.Content.Find.Execute("47÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "55÷9=", 2)

$d.Content.Find.Execute("77÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "77÷7=", 2)
$d.Content.Find.Execute("71÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "26÷8=", 2)
$d.Content.Find.Execute("96÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "88÷2=", 2)
$d.Content.Find.Execute("92÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "59÷7=", 2)
$d.Content.Find.Execute("22÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "52÷3=", 2)

$d.Content.Find.Execute("52÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "59÷6=", 2)
$d.Content.Find.Execute("71÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "27÷8=", 2)
$d.Content.Find.Execute("87÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "15÷6=", 2)
$d.Content.Find.Execute("45÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "57÷6=", 2)
$d.Content.Find.Execute("33÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "66÷6=", 2)
